$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.258.55"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "2.418.33"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.68%  "

$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.106"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("D14").Value = "2.847.63"
$ws.Range("E14").Value = "  -0.07%  "

$ws.Range("D15").Value = "60.161.64"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").Value = "2.393.71"
$ws.Range("E17").Value = "  -0.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.31%  "

$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("E27").Value = "  +2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.96%  "

$ws.Range("D29").Value = "0.0₃0773"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.403"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.20%  "

$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "326.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.94%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0518"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.578"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("E48").Value = "  -1.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "
